# Add 10 more modules to the modules-config-details worksheet.
# Columns: A=Module Name, B=ModuleId, C=difficulty, D=PDF path,
#          E=Author(s), F=Release date, G=count

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("modules-config-details")

$newRows = @(
    @("Shape Memory",        "needyShapeMemory",    1, "modules/Shape Memory.pdf",        "UltraCboy",                   "2017-09-28", 2),
    @("Symbol Cycle",        "SymbolCycleModule",   4, "modules/Symbol Cycle.pdf",         "Timwi",                       "2017-10-05", 1),
    @("Hunting",              "hunting",             3, "modules/Hunting.pdf",              "taggedjc",                    "2017-10-09", 1),
    @("Extended Password",   "ExtendedPassword",    3, "modules/Extended Password.pdf",    "taggedjc, TWGaming",          "2017-10-24", 1),
    @("Curriculum",           "curriculum",          3, "modules/Curriculum.pdf",           "Fixdoll",                     "2017-10-30", 1),
    @("Braille",               "BrailleModule",       4, "modules/Braille.pdf",               "Timwi",                       "2017-10-31", 1),
    @("Mafia",                 "MafiaModule",         4, "modules/Mafia.pdf",                 "MarioXTurn, Ezekiel, Timwi",  "2017-11-04", 1),
    @("Festive Piano Keys",   "FestivePianoKeys",    2, "modules/Festive Piano Keys.pdf",   "Bashly",                      "2017-12-07", 1),
    @("Flags",                 "FlagsModule",         3, "modules/Flags.pdf",                 "Monopoly, Piggered",          "2017-12-24", 1),
    @("Timezone",              "timezone",            3, "modules/Timezone.pdf",              "federan",                     "2017-12-30", 1)
)

$startRow = 146
$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Range($ws.Cells.Item($r, 6), $ws.Cells.Item($r, 6)).NumberFormat = "@"
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

$nextRow = $r

# Mirror Excel's behaviour of scrolling the view and leaving the selection
# on the next empty row below the newly-entered data.
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 136
$win.ScrollColumn = 1
$ws.Cells.Item($nextRow, 1).Select()
